# Apply cryptos.xlsx updates (generated Sat Jun 17 21:39:46 UTC 2023 data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Rows 2-33: update Price (D) and/or Volume(1h) (E) values in place ---
Set-TextValue $ws.Range("D2") "26.492.42"
Set-TextValue $ws.Range("D3") "1.727.33"
Set-TextValue $ws.Range("E3") "  +0.42%  "
Set-TextValue $ws.Range("D4") "0.9998"
Set-TextValue $ws.Range("E4") "  -0.05%  "
Set-TextValue $ws.Range("D5") "245.42"
Set-TextValue $ws.Range("E5") "  +2.71%  "
Set-TextValue $ws.Range("E6") "  -0.07%  "
Set-TextValue $ws.Range("D7") "0.4803"
Set-TextValue $ws.Range("E7") "  +1.86%  "
Set-TextValue $ws.Range("D8") "0.2666"
Set-TextValue $ws.Range("E9") "  +0.26%  "
Set-TextValue $ws.Range("D10") "1.725.65"
Set-TextValue $ws.Range("E10") "  +0.31%  "
Set-TextValue $ws.Range("D11") "0.07149"
Set-TextValue $ws.Range("E11") "  +1.10%  "
Set-TextValue $ws.Range("E12") "  +2.29%  "
Set-TextValue $ws.Range("D13") "0.6153"
Set-TextValue $ws.Range("E13") "  +4.16%  "
Set-TextValue $ws.Range("D14") "4.519"
Set-TextValue $ws.Range("E14") "  +2.77%  "
Set-TextValue $ws.Range("D15") "77.12"
Set-TextValue $ws.Range("E15") "  +1.27%  "
Set-TextValue $ws.Range("E16") "  -0.07%  "
Set-TextValue $ws.Range("D17") "26.495.32"
Set-TextValue $ws.Range("E17") "  +0.66%  "
Set-TextValue $ws.Range("D19") "0.000006931"
Set-TextValue $ws.Range("E19") "  +2.22%  "
Set-TextValue $ws.Range("E20") "  +0.69%  "
Set-TextValue $ws.Range("D21") "1.946.82"
Set-TextValue $ws.Range("E21") "  +0.35%  "
Set-TextValue $ws.Range("D22") "4.526"
Set-TextValue $ws.Range("E22") "  -0.54%  "
Set-TextValue $ws.Range("D23") "8.949"
Set-TextValue $ws.Range("E23") "  +2.61%  "
Set-TextValue $ws.Range("D24") "5.275"
Set-TextValue $ws.Range("E24") "  -0.95%  "
Set-TextValue $ws.Range("D25") "136.73"
Set-TextValue $ws.Range("E25") "  +0.88%  "
Set-TextValue $ws.Range("D26") "15.35"
Set-TextValue $ws.Range("E26") "  +0.86%  "
Set-TextValue $ws.Range("D27") "1.797"
Set-TextValue $ws.Range("E27") "  +2.28%  "
Set-TextValue $ws.Range("D28") "1.404"
Set-TextValue $ws.Range("E28") "  -0.43%  "
Set-TextValue $ws.Range("D29") "106.86"
Set-TextValue $ws.Range("E29") "  -1.37%  "
Set-TextValue $ws.Range("D30") "3.974"
Set-TextValue $ws.Range("E30") "  -0.72%  "
Set-TextValue $ws.Range("D31") "0.08030"
Set-TextValue $ws.Range("E31") "  +3.83%  "
Set-TextValue $ws.Range("D32") "3.707"
Set-TextValue $ws.Range("E32") "  +0.59%  "
Set-TextValue $ws.Range("D33") "0.04563"
Set-TextValue $ws.Range("E33") "  +2.29%  "

# --- Rows 34-51: a new "Frax" entry is inserted at row 34, shifting the
# remaining coins (HuobiToken..NEARProtocol) down by one row, which pushes
# the former last row (Decentraland) off the bottom of the table. The rank
# numbers in column A stay fixed per row, so only columns B-E are rewritten.
$ws.Range("B34").Value = "Frax"
$ws.Range("C34").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
Set-TextValue $ws.Range("D34") "0.9998"
Set-TextValue $ws.Range("E34") "  -0.10%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue $ws.Range("D35") "2.615"
Set-TextValue $ws.Range("E35") "  -0.02%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue $ws.Range("D36") "0.6353"
Set-TextValue $ws.Range("E36") "  +2.70%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D37") "0.9924"
Set-TextValue $ws.Range("E37") "  +1.92%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D38") "0.9307"
Set-TextValue $ws.Range("E38") "  +0.91%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D39") "2.087"
Set-TextValue $ws.Range("E39") "  +9.96%  "
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D40") "2.425"
Set-TextValue $ws.Range("E40") "  +0.35%  "
$ws.Range("B41").Value = "Quant"
$ws.Range("C41").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D41") "105.31"
Set-TextValue $ws.Range("E41") "  -8.17%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue $ws.Range("D42") "1.006"
Set-TextValue $ws.Range("E42") "  +0.41%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D43") "0.01503"
Set-TextValue $ws.Range("E43") "  +1.87%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D44") "5.583"
Set-TextValue $ws.Range("E44") "  +4.28%  "
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D45") "0.3899"
Set-TextValue $ws.Range("E45") "  +2.29%  "
$ws.Range("B46").Value = "Aptos"
$ws.Range("C46").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D46") "6.936"
Set-TextValue $ws.Range("E46") "  +10.75%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D47") "0.1183"
Set-TextValue $ws.Range("E47") "  +1.62%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D48") "0.05335"
Set-TextValue $ws.Range("E48") "  +0.90%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue $ws.Range("D49") "30.91"
Set-TextValue $ws.Range("E49") "  +1.10%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.877"
Set-TextValue $ws.Range("E50") "  +2.27%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D51") "1.267"
Set-TextValue $ws.Range("E51") "  +4.40%  "

Write-Host "cryptos.xlsx updated"
